# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-22) is re-sorted in ascending order by
# the "Periodo Mora" column (E), and the "Salario Basico" column (G) for
# every one of those rows is updated to the new base salary value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sort the data rows (B16:J22) ascending by the "Periodo Mora" column (E).
$dataRange = $ws.Range("B16:J22")
$sortKey = $ws.Range("E16:E22")
$dataRange.Sort($sortKey, 1)

# Refresh "Salario Basico" (column G) for all rows in the table with the
# updated value.
$ws.Range("G16:G22").Value = 828116
